# Weekly update: insert a new price record as row 10, pushing the
# existing rows (old 10..61) down by one (new 11..62).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 10; this shifts rows 10-61 down to 11-62
# and the sheet's used-range grows from A1:R61 to A1:R62 automatically.
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with the new weekly record.
# All the "constant" columns match every other row in this sheet.
$ws.Range("A10").Value = 10
$ws.Range("B10").Value = "Vega Modelo de Temuco"
$ws.Range("C10").Value = "La Araucanía"
$ws.Range("D10").Value = "2023-02-02"
$ws.Range("E10").Value = 9
$ws.Range("F10").Value = 100112042
$ws.Range("G10").Value = "Locoto"
$ws.Range("H10").Value = "Sin especificar"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 140
$ws.Range("K10").Value = 2500
$ws.Range("L10").Value = 2500
$ws.Range("M10").Value = 2500
$ws.Range("N10").Value = "$/kilo"
$ws.Range("O10").Value = "Región de Arica y Parinacota"
$ws.Range("P10").Value = 2500
$ws.Range("Q10").Value = 1
$ws.Range("R10").Value = "Hortaliza"
